{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst sec = sections.items[0];\nconst footer = sec.getFooter(\"Primary\");\nconst body = footer.getRange();\nbody.insertText(\"X\", \"End\");\nawait context.sync();\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n$sec = $d.Sections(1)\n$header = $sec.Headers(1)\n$r = $header.Range\n$x = $r.InsertXML($hfrag)\n"}
